$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2696.5686
$ws.Range("J17").Value = 2852.6667
$ws.Range("L17").Value = 8558.000100000001
$ws.Range("N17").Value = -8894.000100000001

$ws.Range("H33").Value = 83.69231000000001
$ws.Range("I33").Value = 88
$ws.Range("J33").Value = 60
$ws.Range("K33").Value = 88
$ws.Range("L33").Value = 60
$ws.Range("M33").Value = 141
$ws.Range("N33").Value = -518

$ws.Range("H43").Value = 1194.4375
$ws.Range("I43").Value = 979.375
$ws.Range("J43").Value = 1409.5
$ws.Range("K43").Value = 979.375
$ws.Range("L43").Value = 1409.5
$ws.Range("M43").Value = -910.375
$ws.Range("N43").Value = -1547.5

$ws.Range("H64").Value = 3810.9
$ws.Range("I64").Value = 2744
$ws.Range("J64").Value = 4385.385
$ws.Range("K64").Value = 2744
$ws.Range("L64").Value = 4385.385
$ws.Range("M64").Value = -2496
$ws.Range("N64").Value = -4881.385

$ws.Range("H67").Value = 3810.9
$ws.Range("I67").Value = 2744
$ws.Range("J67").Value = 4385.385
$ws.Range("K67").Value = 2744
$ws.Range("L67").Value = 4385.385
$ws.Range("M67").Value = -1886
$ws.Range("N67").Value = -6101.385

$ws.Range("H86").Value = 2689.05
$ws.Range("I86").Value = 2009.1333
$ws.Range("J86").Value = 4728.8
$ws.Range("K86").Value = 2009.1333
$ws.Range("L86").Value = 4728.8
$ws.Range("M86").Value = -886.1333
$ws.Range("N86").Value = -6974.8

$ws.Range("H89").Value = 2689.05
$ws.Range("I89").Value = 2009.1333
$ws.Range("J89").Value = 4728.8
$ws.Range("K89").Value = 10045.6665
$ws.Range("L89").Value = 23644
$ws.Range("M89").Value = -4429.666499999999
$ws.Range("N89").Value = -34876

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 67638.53
$ws.Range("I2").Value = 84189.914
$ws.Range("J2").Value = 1433
$ws.Range("K2").Value = 84189.914
$ws.Range("L2").Value = 1433
$ws.Range("M2").Value = -84076.914
$ws.Range("N2").Value = -1659

$ws.Range("H74").Value = 2939.8
$ws.Range("I74").Value = 782.18604
$ws.Range("J74").Value = 16193.714
$ws.Range("K74").Value = 782.18604
$ws.Range("L74").Value = 16193.714
$ws.Range("M74").Value = 91.81395999999995
$ws.Range("N74").Value = -17941.714

$ws.Range("H77").Value = 2939.8
$ws.Range("I77").Value = 782.18604
$ws.Range("J77").Value = 16193.714
$ws.Range("K77").Value = 3910.9302
$ws.Range("L77").Value = 80968.57000000001
$ws.Range("M77").Value = 457.0697999999998
$ws.Range("N77").Value = -89704.57000000001

$ws.Range("H102").Value = 1188
$ws.Range("I102").Value = 970.3570999999999
$ws.Range("J102").Value = 1797.4
$ws.Range("K102").Value = 970.3570999999999
$ws.Range("L102").Value = 1797.4
$ws.Range("M102").Value = 651.6429000000001
$ws.Range("N102").Value = -5041.4

$ws.Range("H116").Value = 67638.53
$ws.Range("I116").Value = 84189.914
$ws.Range("J116").Value = 1433
$ws.Range("K116").Value = 84189.914
$ws.Range("L116").Value = 1433
$ws.Range("M116").Value = -81895.914
$ws.Range("N116").Value = -6021

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 67638.53
$ws.Range("I3").Value = 84189.914
$ws.Range("J3").Value = 1433
$ws.Range("K3").Value = 84189.914
$ws.Range("L3").Value = 1433
$ws.Range("M3").Value = -84075.914
$ws.Range("N3").Value = -1661

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4200.6875
$ws.Range("I62").Value = 5428.0557
$ws.Range("J62").Value = 2622.6428
$ws.Range("K62").Value = 5428.0557
$ws.Range("L62").Value = 2622.6428
$ws.Range("M62").Value = -4804.0557
$ws.Range("N62").Value = -3870.6428

$ws.Range("H65").Value = 4200.6875
$ws.Range("I65").Value = 5428.0557
$ws.Range("J65").Value = 2622.6428
$ws.Range("K65").Value = 27140.2785
$ws.Range("L65").Value = 13113.214
$ws.Range("M65").Value = -24020.2785
$ws.Range("N65").Value = -19353.214

$ws.Range("H100").Value = 23000
$ws.Range("J100").Value = 23000
$ws.Range("L100").Value = 23000
$ws.Range("N100").Value = -25164

$ws.Range("H134").Value = 1385.4565
$ws.Range("I134").Value = 1298.3823
$ws.Range("J134").Value = 1632.1666
$ws.Range("K134").Value = 3895.1469
$ws.Range("L134").Value = 4896.4998
$ws.Range("M134").Value = -1360.1469
$ws.Range("N134").Value = -9966.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1004.4545
$ws.Range("I34").Value = 289
$ws.Range("J34").Value = 1413.2858
$ws.Range("K34").Value = 867
$ws.Range("L34").Value = 4239.857400000001
$ws.Range("M34").Value = -783
$ws.Range("N34").Value = -4407.857400000001

$ws.Range("H39").Value = 3371.2856
$ws.Range("J39").Value = 3833.3333
$ws.Range("L39").Value = 11499.9999
$ws.Range("N39").Value = -12087.9999

$ws.Range("H55").Value = 2394.1667
$ws.Range("I55").Value = 1354.4445
$ws.Range("J55").Value = 3433.889
$ws.Range("K55").Value = 4063.3335
$ws.Range("L55").Value = 10301.667
$ws.Range("M55").Value = -3886.3335
$ws.Range("N55").Value = -10655.667

$ws.Range("H132").Value = 1132.7693
$ws.Range("I132").Value = 948
$ws.Range("K132").Value = 8532
$ws.Range("M132").Value = -6002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4952.6665
$ws.Range("J70").Value = 5068.1665
$ws.Range("L70").Value = 5068.1665
$ws.Range("N70").Value = -5608.1665

$ws.Range("H73").Value = 4952.6665
$ws.Range("J73").Value = 5068.1665
$ws.Range("L73").Value = 5068.1665
$ws.Range("N73").Value = -6940.1665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1373.0667
$ws.Range("I61").Value = 1161.375
$ws.Range("J61").Value = 1615
$ws.Range("K61").Value = 1161.375
$ws.Range("L61").Value = 1615
$ws.Range("M61").Value = -959.375
$ws.Range("N61").Value = -2019

$ws.Range("H93").Value = 1218.2609
$ws.Range("I93").Value = 1407.3334
$ws.Range("J93").Value = 1012
$ws.Range("K93").Value = 1407.3334
$ws.Range("L93").Value = 1012
$ws.Range("M93").Value = -159.3334
$ws.Range("N93").Value = -3508

$ws.Range("H113").Value = 1373.0667
$ws.Range("I113").Value = 1161.375
$ws.Range("J113").Value = 1615
$ws.Range("K113").Value = 1161.375
$ws.Range("L113").Value = 1615
$ws.Range("M113").Value = 1008.625
$ws.Range("N113").Value = -5955

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2492.0598
$ws.Range("I132").Value = 2661.7322
$ws.Range("K132").Value = 7985.196599999999
$ws.Range("M132").Value = -5455.196599999999
